$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-21 03:48:40'
$ws.Range('H2').NumberFormat = '@'
$ws.Range('H2').Value = '58%'
$ws.Range('E3').Value = '2026-02-21 03:48:43'
$ws.Range('M3').Value = '-1.0 °C 3:22 TU'
$ws.Range('O3').Value = '-2.1 °C'
$ws.Range('E4').Value = '2026-02-21 03:48:45'
$ws.Range('N4').Value = '1.2 °C 3:29 TU'
$ws.Range('O4').Value = '2.3 °C'
$ws.Range('E5').Value = '2026-02-21 03:48:48'
$ws.Range('G5').Value = '131 cm'
$ws.Range('H5').NumberFormat = '@'
$ws.Range('H5').Value = '68%'
$ws.Range('M5').Value = '0.4 °C 3:28 TU'
$ws.Range('O5').Value = '-1.5 °C'
$ws.Range('E6').Value = '2026-02-21 03:48:50'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '90%'
$ws.Range('N6').Value = '3.7 °C 3:17 TU'
$ws.Range('O6').Value = '5.1 °C'
$ws.Range('E7').Value = '2026-02-21 03:48:53'
$ws.Range('N7').Value = '10.7 °C 3:03 TU'
$ws.Range('O7').Value = '11.6 °C'
$ws.Range('E8').Value = '2026-02-21 03:48:55'
$ws.Range('N8').Value = '7.0 °C 3:27 TU'
$ws.Range('O8').Value = '7.8 °C'
$ws.Range('E9').Value = '2026-02-21 03:48:57'
$ws.Range('N9').Value = '11.2 °C 3:19 TU'
$ws.Range('O9').Value = '11.8 °C'
$ws.Range('E10').Value = '2026-02-21 03:49:00'
$ws.Range('N10').Value = '1.4 °C 3:24 TU'
$ws.Range('O10').Value = '2.4 °C'
$ws.Range('E11').Value = '2026-02-21 03:49:03'
$ws.Range('O11').Value = '4.5 °C'
$ws.Range('E12').Value = '2026-02-21 03:49:06'
$ws.Range('N12').Value = '10.0 °C 3:23 TU'
$ws.Range('O12').Value = '12.1 °C'
$ws.Range('E13').Value = '2026-02-21 03:49:08'
$ws.Range('J13').Value = '1035.2 hPa'
$ws.Range('N13').Value = '-4.5 °C 3:19 TU'
$ws.Range('O13').Value = '-2.8 °C'
$ws.Range('E14').Value = '2026-02-21 03:49:11'
$ws.Range('H14').NumberFormat = '@'
$ws.Range('H14').Value = '71%'
$ws.Range('N14').Value = '7.3 °C 3:19 TU'
$ws.Range('O14').Value = '9.4 °C'
$ws.Range('E15').Value = '2026-02-21 03:49:14'
$ws.Range('N15').Value = '11.5 °C 3:20 TU'
$ws.Range('O15').Value = '11.9 °C'
$ws.Range('E16').Value = '2026-02-21 03:49:16'
$ws.Range('H16').NumberFormat = '@'
$ws.Range('H16').Value = '36%'
$ws.Range('E17').Value = '2026-02-21 03:49:18'
$ws.Range('E18').Value = '2026-02-21 03:49:21'
$ws.Range('J18').Value = '1029.2 hPa'
$ws.Range('N18').Value = '0.3 °C 3:29 TU'
$ws.Range('O18').Value = '1.2 °C'
$ws.Range('E19').Value = '2026-02-21 03:49:24'
$ws.Range('H19').NumberFormat = '@'
$ws.Range('H19').Value = '87%'
$ws.Range('N19').Value = '2.1 °C 3:16 TU'
$ws.Range('O19').Value = '2.9 °C'
$ws.Range('E20').Value = '2026-02-21 03:49:26'
$ws.Range('E21').Value = '2026-02-21 03:49:28'
$ws.Range('J21').Value = '1032.6 hPa'
$ws.Range('O21').Value = '1.1 °C'
$ws.Range('E22').Value = '2026-02-21 03:49:31'
$ws.Range('H22').NumberFormat = '@'
$ws.Range('H22').Value = '43%'
$ws.Range('E23').Value = '2026-02-21 03:49:34'
$ws.Range('H23').NumberFormat = '@'
$ws.Range('H23').Value = '39%'
$ws.Range('O23').Value = '-0.2 °C'
$ws.Range('E24').Value = '2026-02-21 03:49:36'
$ws.Range('J24').Value = '1030.9 hPa'
$ws.Range('N24').Value = '0.3 °C 3:23 TU'
$ws.Range('O24').Value = '2.5 °C'
$ws.Range('E25').Value = '2026-02-21 03:49:39'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '47%'
$ws.Range('E26').Value = '2026-02-21 03:49:42'
$ws.Range('H26').NumberFormat = '@'
$ws.Range('H26').Value = '38%'
$ws.Range('E27').Value = '2026-02-21 03:49:45'
$ws.Range('N27').Value = '0.5 °C 3:15 TU'
$ws.Range('O27').Value = '0.8 °C'
$ws.Range('E28').Value = '2026-02-21 03:49:47'
$ws.Range('L28').Value = '8.6 km/h - 232º 3:01 TU'
$ws.Range('E29').Value = '2026-02-21 03:49:50'
$ws.Range('E30').Value = '2026-02-21 03:49:53'
$ws.Range('J30').Value = '1027.9 hPa'
$ws.Range('O30').Value = '10.2 °C'
$ws.Range('E31').Value = '2026-02-21 03:49:56'
$ws.Range('J31').Value = '1026.7 hPa'
$ws.Range('N31').Value = '9.7 °C 3:03 TU'
$ws.Range('O31').Value = '10.1 °C'
$ws.Range('E32').Value = '2026-02-21 03:49:58'
$ws.Range('N32').Value = '0.7 °C 3:25 TU'
$ws.Range('O32').Value = '1.4 °C'
$ws.Range('E33').Value = '2026-02-21 03:50:01'
$ws.Range('J33').Value = '1033.0 hPa'
$ws.Range('N33').Value = '-1.2 °C 3:29 TU'
$ws.Range('O33').Value = '-0.1 °C'
$ws.Range('E34').Value = '2026-02-21 03:50:04'
$ws.Range('H34').NumberFormat = '@'
$ws.Range('H34').Value = '39%'
$ws.Range('L34').Value = '45.4 km/h - 18º 3:06 TU'
$ws.Range('M34').Value = '3.8 °C 3:17 TU'
$ws.Range('O34').Value = '3.1 °C'
$ws.Range('E35').Value = '2026-02-21 03:50:06'
$ws.Range('N35').Value = '3.1 °C 3:27 TU'
$ws.Range('E36').Value = '2026-02-21 03:50:09'
$ws.Range('J36').Value = '1027.7 hPa'
$ws.Range('N36').Value = '12.1 °C 3:24 TU'
$ws.Range('O36').Value = '12.5 °C'
$ws.Range('E37').Value = '2026-02-21 03:50:12'
$ws.Range('N37').Value = '-1.7 °C 3:00 TU'
$ws.Range('O37').Value = '-1.3 °C'
$ws.Range('E38').Value = '2026-02-21 03:50:14'
$ws.Range('H38').NumberFormat = '@'
$ws.Range('H38').Value = '92%'
$ws.Range('N38').Value = '3.3 °C 3:08 TU'
$ws.Range('O38').Value = '4.4 °C'
$ws.Range('E39').Value = '2026-02-21 03:50:17'
$ws.Range('H39').NumberFormat = '@'
$ws.Range('H39').Value = '41%'
$ws.Range('O39').Value = '0.0 °C'
$ws.Range('E40').Value = '2026-02-21 03:50:19'
$ws.Range('H40').NumberFormat = '@'
$ws.Range('H40').Value = '75%'
$ws.Range('J40').Value = '1032.3 hPa'
$ws.Range('M40').Value = '4.6 °C 3:29 TU'
$ws.Range('O40').Value = '2.3 °C'
$ws.Range('E41').Value = '2026-02-21 03:50:22'
$ws.Range('H41').NumberFormat = '@'
$ws.Range('H41').Value = '64%'
$ws.Range('N41').Value = '6.2 °C 3:29 TU'
$ws.Range('O41').Value = '9.6 °C'
$ws.Range('E42').Value = '2026-02-21 03:50:25'
$ws.Range('O42').Value = '8.2 °C'
$ws.Range('E43').Value = '2026-02-21 03:50:27'
$ws.Range('N43').Value = '-0.3 °C 3:29 TU'
$ws.Range('O43').Value = '0.7 °C'
$ws.Range('E44').Value = '2026-02-21 03:50:30'
$ws.Range('H44').NumberFormat = '@'
$ws.Range('H44').Value = '47%'
$ws.Range('M44').Value = '1.3 °C 3:29 TU'
$ws.Range('O44').Value = '0.4 °C'
$ws.Range('E45').Value = '2026-02-21 03:50:32'
$ws.Range('H45').NumberFormat = '@'
$ws.Range('H45').Value = '87%'
$ws.Range('J45').Value = '1035.2 hPa'
$ws.Range('N45').Value = '-1.3 °C 3:09 TU'
$ws.Range('O45').Value = '0.7 °C'
$ws.Range('E46').Value = '2026-02-21 03:50:34'
$ws.Range('H46').NumberFormat = '@'
$ws.Range('H46').Value = '82%'
$ws.Range('N46').Value = '4.7 °C 3:00 TU'
$ws.Range('O46').Value = '5.9 °C'
